$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column: several values look like plain numbers (e.g. "554.39")
# while others use "." as a thousands separator (e.g. "60.004.56") or contain
# non-numeric characters. To keep every D-column value as literal text (matching
# the source inline strings) force each target cell to Text format first.
$dCells = @("D2", "D3", "D5", "D6", "D10", "D12", "D13", "D14", "D15", "D17", "D18", "D20", "D24", "D25", "D29", "D30", "D32", "D34", "D38", "D39", "D41", "D44", "D48")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "60.004.56"
$ws.Range("D3").Value = "2.409.79"
$ws.Range("D5").Value = "554.39"
$ws.Range("D6").Value = "136.07"
$ws.Range("D10").Value = "5.61"
$ws.Range("D12").Value = "0.351"
$ws.Range("D13").Value = "24.67"
$ws.Range("D14").Value = "2.839.76"
$ws.Range("D15").Value = "59.894.26"
$ws.Range("D17").Value = "2.407.22"
$ws.Range("D18").Value = "11.17"
$ws.Range("D20").Value = "327.49"
$ws.Range("D24").Value = "0.180"
$ws.Range("D25").Value = "8.60"
$ws.Range("D29").Value = "0.0₃0768"
$ws.Range("D30").Value = "170.78"
$ws.Range("D32").Value = "1.11"
$ws.Range("D34").Value = "18.41"
$ws.Range("D38").Value = "4.20"
$ws.Range("D39").Value = "321.91"
$ws.Range("D41").Value = "146.24"
$ws.Range("D44").Value = "19.89"
$ws.Range("D48").Value = "11.05"

# Restore the default (Normal) style now that the text is stored, so no stray
# number-format style lingers on these cells.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Volume(1h) (E) column updates.
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  +3.79%  "
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("E24").Value = "  +4.41%  "
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +4.74%  "
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("E32").Value = "  +9.14%  "
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("E41").Value = "  +5.84%  "
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("E51").Value = "  -2.04%  "
